# Auto-generated edit script for cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "25.632.70"
$ws.Cells.Item(2, 4).ClearFormats()
$ws.Cells.Item(2, 5).Value = "  -4.23%  "

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.812.72"
$ws.Cells.Item(3, 4).ClearFormats()
$ws.Cells.Item(3, 5).Value = "  -3.05%  "

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.9990"
$ws.Cells.Item(4, 4).ClearFormats()
$ws.Cells.Item(4, 5).Value = "  -0.21%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "275.92"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = "  -8.04%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.9985"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = "  -0.25%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.5045"
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(7, 5).Value = "  -5.26%  "

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3437"
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).Value = "  -7.79%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "44.00"
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).Value = "  -2.97%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.06850"
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).Value = "  -4.11%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "19.67"
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).Value = "  -8.42%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.8018"
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).Value = "  -9.53%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.07810"
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).Value = "  -4.15%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "1.813.15"
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).Value = "  -3.53%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "5.047"
$ws.Cells.Item(15, 4).ClearFormats()
$ws.Cells.Item(15, 5).Value = "  -4.61%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "87.45"
$ws.Cells.Item(16, 4).ClearFormats()
$ws.Cells.Item(16, 5).Value = "  -5.39%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.9983"
$ws.Cells.Item(17, 4).ClearFormats()
$ws.Cells.Item(17, 5).Value = "  -0.28%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "14.02"
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(18, 5).Value = "  -5.51%  "

# Row 19
$ws.Cells.Item(19, 2).Value = "ShibaInu"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.000008023"
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).Value = "  -5.57%  "

# Row 20
$ws.Cells.Item(20, 2).Value = "Dai"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.9988"
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).Value = "  -0.19%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "25.680.07"
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(21, 5).Value = "  -4.17%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "4.722"
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).Value = "  -5.07%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "9.928"
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(23, 5).Value = "  -6.56%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "6.115"
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).Value = "  -4.12%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.282"
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value = "  -0.13%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "142.32"
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 5).Value = "  -2.29%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  -4.61%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "17.08"
$ws.Cells.Item(28, 4).ClearFormats()
$ws.Cells.Item(28, 5).Value = "  -5.18%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "108.76"
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).Value = "  -4.19%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "4.284"
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(30, 5).Value = "  -8.83%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "4.228"
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).Value = "  -8.57%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.08724"
$ws.Cells.Item(32, 4).ClearFormats()
$ws.Cells.Item(32, 5).Value = "  -4.08%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.04801"
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).Value = "  -4.40%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.139"
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = "  -2.74%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.7215"
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(35, 5).Value = "  -10.79%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "2.837"
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(36, 5).Value = "  -3.73%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "3.132"
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).Value = "  -1.56%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "2.354"
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(38, 5).Value = "  -11.21%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.01836"
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).Value = "  -5.38%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.5075"
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).Value = "  -16.91%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.9393"
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).Value = "  -11.98%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "115.84"
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).Value = "  +0.42%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "6.171"
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).Value = "  -4.58%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "7.851"
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).Value = "  -10.05%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.9975"
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).Value = "  -0.35%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.1354"
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).Value = "  -9.13%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.4428"
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(47, 5).Value = "  -16.55%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "EnergySwap"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "9.271"
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(48, 5).Value = "  -6.79%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "Elrond"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "36.14"
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).Value = "  -3.13%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.05895"
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).Value = "  -2.64%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.479"
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).Value = "  -10.39%  "
